$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.ActiveWindow.Zoom = 115
